$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.831.79'
$ws.Range('E2').Value = '  +2.63%  '
$ws.Range('D3').Value = '1.881.91'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.006'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.56%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '323.94'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.52%  '
$ws.Range('E6').Value = '  +0.49%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4669'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.47%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3932'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.72%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07923'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.72%  '
$ws.Range('E10').Value = '  +2.32%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '22.40'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.87%  '
$ws.Range('D12').Value = '1.823.17'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.747'
$ws.Range('D13').ClearFormats()
$ws.Range('E14').Value = '  +1.85%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.06985'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.94%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '88.73'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.37%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.006'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('E18').Value = '  +1.18%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '16.97'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.80%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.005'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('D21').Value = '28.854.47'
$ws.Range('E21').Value = '  +2.59%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.349'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.10'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.84%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.126'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.53%  '
$ws.Range('D25').Value = '2.150.24'
$ws.Range('E25').Value = '  +4.38%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '153.23'
$ws.Range('D26').ClearFormats()
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.38'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.83%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '5.809'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.22%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.003'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.49%  '
$ws.Range('E30').Value = '  +2.47%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.09396'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.41%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.9387'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.314'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.17%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.359'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.97%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.346'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.05920'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.41%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02128'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.02%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.162'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.27%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '7.897'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.99%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.5734'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.52%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.1798'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.43%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '10.01'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.85%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.07315'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.98%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '11.89'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.46%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.176'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.58%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5362'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.60%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.848'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.73%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '114.11'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.24%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.089'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -5.84%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.376'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.03%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.004'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.55%  '
